# "modify progress add some music"
# Append three new progress-log rows (83-85) to Sheet1, right after the
# existing last row (82), mirroring the same A/B/C/D layout used by the
# rest of the table (A=date, B=target/work done, C=risk, D=duration in
# hours).
#
# Column A holds date-looking text (e.g. "2012.11.25") that must stay as
# literal text, not get auto-converted into a serial date number. The
# trick: assign it through .Formula with a leading apostrophe (forces
# text), then restore the normal cell formatting by pasting-special the
# (already plain) format from the row above, so the new cells end up
# looking just like every other date cell in the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 83 ----------------------------------------------------------
$ws.Range("A83").Formula = "'2012.11.25"
$ws.Range("A82").Copy()
$ws.Range("A83").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B83").Value = "代码结构重构，ipad同步"
$ws.Range("D83").Value = 3

# ---- Row 84 ----------------------------------------------------------
$ws.Range("A84").Formula = "'2012.12.2"
$ws.Range("A82").Copy()
$ws.Range("A84").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B84").Value = "结束界面，字体，特效，"
$ws.Range("C84").Value = "刘云鹏，赵羽佳，刘晋 集中开发"
$ws.Range("D84").Value = 6

# ---- Row 85 ----------------------------------------------------------
$ws.Range("B85").Value = "角色进度条修改"
$ws.Range("D85").Value = 2

# Match the author's final cursor position recorded in the saved file.
$ws.Range("C88").Select()
